# Commit: "#5: property aircraft done"
#
# The 建物 (building) sheet and 汽車 (car) sheet had their
# "property_category" column incorrectly populated with the value "land"
# (carried over from the 土地 sheet). Fix the category labels so each
# sheet reports its own property type:
#   - 建物 (Buildings) sheet: property_category -> "building"
#   - 汽車 (Cars) sheet:      property_category -> "car"

$wb = $excel.ActiveWorkbook

# --- 建物 (Buildings) sheet ---------------------------------------------
$wsBuilding = $wb.Worksheets.Item("建物")

# Column I holds "property_category" (header in I1); rows 2-7 are data.
for ($row = 2; $row -le 7; $row++) {
    $wsBuilding.Cells.Item($row, 9).Value = "building"
}

# --- 汽車 (Cars) sheet ---------------------------------------------------
$wsCar = $wb.Worksheets.Item("汽車")

# Column H holds "property_category" (header in H1); row 2 is the only data row.
$wsCar.Cells.Item(2, 8).Value = "car"
